$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''25.792.93'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -0.82%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.628.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -0.84%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.14%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''215.37'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +0.02%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''0.5057'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.14%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.50%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.06415'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  +0.72%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''19.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -2.17%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.07790'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  +0.52%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''4.258'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -1.04%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''1.628.67'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -0.79%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '''1.852.42'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.92%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = '''0.5584'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  +1.81%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''63.00'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  -2.27%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''0.0₅7556'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -2.65%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''25.799.54'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E19').Value = '''  +0.09%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''194.10'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -1.94%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''4.324'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Value = '''9.828'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -1.58%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''5.996'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -2.43%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = '''  +0.08%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''1.794'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -5.29%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = '''140.78'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -1.32%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = '''0.1269'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +0.85%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''6.732'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -2.20%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''15.40'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -1.75%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = '''  -0.43%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = '''0.04866'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -0.78%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''3.283'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -0.18%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''3.201'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -0.46%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = '''1.555'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -0.32%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '''  -0.12%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''0.8947'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -2.81%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''2.573'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  +0.14%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''1.128.03'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  +1.81%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = '''0.5464'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -1.84%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = '''0.01560'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -0.71%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''0.9951'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -0.70%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''5.546'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  -1.30%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''0.7976'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -0.95%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''97.23'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -1.61%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''1.778.61'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -0.14%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = '''  -4.21%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = '''0.4442'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -2.03%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''55.26'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -0.23%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '''  -2.71%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''7.686'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Value = '''1.001'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  -0.18%  '
$ws.Range('E51').Style = 'Normal'
